$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("final_figures/2009/F1_P21_Barton & Sutcliffe_2009_HR_Overcoming dysfunctional mometum - organizational safety as a social achievement.png", "Yes", "Yes", "No", "No", "No", "Yes", "No"),
    @("final_figures/2009/F1_P18_Bunderson & Thompson_2009_ASQ_The Call of the Wild.png", "Yes", "Yes", "No", "No", "No", "No", "Yes"),
    @("final_figures/2009/F1_P13_Elsbach_2009_HR_Identity affirmation through signature style - a study of toy car designers.png", "Yes", "Yes", "No", "No", "Yes", "Yes", "No"),
    @("final_figures/2009/F1_P7_Gordon, Clegg & Kornberger_2009_OrgStudies_Embedded Ethics - Discourse and Power in the New south Wales Police Service.png", "Yes", "Yes", "No", "No", "Yes", "Yes", "No"),
    @("final_figures/2009/F2_P2_Gulati & Puranam_2009_OrgSci_Renewal Through Reorganization.png", "Yes", "No", "No", "No", "No", "Yes", "No"),
    @("final_figures/2009/F2_P25_Levina & Orlikowski_2009_AMJ_Understanding Shifting Power Relations within and across Organizations a Critical Genre Analysis.png", "Yes", "Yes", "No", "No", "No", "Yes", "Yes"),
    @("final_figures/2009/F1_P5_Rerup_2009_OrgSci_Attentional Triangulation.png", "Yes", "Yes", "No", "Yes", "Yes", "Yes", "No"),
    @("final_figures/2009/F1_P7_Sonenshein_2009_OrgSci_Emergence of Ethical Issues During Strategic Change Implementation.png", "Yes", "Yes", "No", "Yes", "No", "No", "Yes"),
    @("final_figures/2009/F1_P7_Kreiner, Hollensbe & Sheep_2009_AMJ_Balancing Borders and Bridges Negotiating the Work-Home Interface Via Boundary Work Tactics.png", "Yes", "Yes", "No", "Yes", "No", "Yes", "No"),
    @("final_figures/2009/F1_P10_Levina & Orlikowski_2009_AMJ_Understanding Shifting Power Relations within and across Organizations a Critical Genre Analysis.png", "Yes", "No", "No", "No", "Yes", "No", "Yes")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
